$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The DOB column (C2:C6) held the literal text "2025-06-12" (shared string) for
# every employee row. Replace it with a real date value (2025-06-13, Excel
# serial 45821) formatted with the built-in short-date number format.
$dobRange = $ws.Range("C2:C6")
$dobRange.Value = 45821

# Apply the built-in "m/d/yy"-style date format (numFmtId 14) to the first
# cell, then propagate that exact style to the rest of the column via a
# format-only paste so every cell in C2:C6 shares a single cell style.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C2:C6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column C is now wide enough to show the date without truncation.
$ws.Columns.Item(3).ColumnWidth = 8.7109375

# Restore the last active selection used when the file was saved.
$ws.Range("G20").Select() | Out-Null
